$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MAE/MSE metric values (rows reporting MSE/MAE for various light-sampling configs)
$ws.Range("B10").Value = 0.004472
$ws.Range("C10").Value = 0.004812
$ws.Range("D10").Value = 0.000878
$ws.Range("E10").Value = 0.001307
$ws.Range("F10").Value = 0.011591
$ws.Range("G10").Value = 0.010984
$ws.Range("H10").Value = 0.025195
$ws.Range("I10").Value = 0.166043

$ws.Range("B11").Value = 0.030145
$ws.Range("C11").Value = 0.019896
$ws.Range("D11").Value = 0.005852
$ws.Range("E11").Value = 0.006712
$ws.Range("F11").Value = 0.023533
$ws.Range("G11").Value = 0.021996
$ws.Range("H11").Value = 0.01107
$ws.Range("I11").Value = 0.014212

$ws.Range("B14").Value = 0.00629
$ws.Range("C14").Value = 0.007987
$ws.Range("D14").Value = 0.009716
$ws.Range("E14").Value = 0.002195
$ws.Range("F14").Value = 0.008631
$ws.Range("G14").Value = 0.111707
$ws.Range("H14").Value = 0.008625
$ws.Range("I14").Value = 0.012742

$ws.Range("B15").Value = 0.031135
$ws.Range("C15").Value = 0.020677
$ws.Range("D15").Value = 0.005958
$ws.Range("E15").Value = 0.006817
$ws.Range("F15").Value = 0.023898
$ws.Range("G15").Value = 0.022954
$ws.Range("H15").Value = 0.011056
$ws.Range("I15").Value = 0.01394

$ws.Range("B18").Value = 0.008836
$ws.Range("C18").Value = 0.012676
$ws.Range("D18").Value = 0.002577
$ws.Range("E18").Value = 0.001511
$ws.Range("F18").Value = 0.084108
$ws.Range("G18").Value = 0.038126
$ws.Range("H18").Value = 0.006793
$ws.Range("I18").Value = 0.053589

$ws.Range("B19").Value = 0.032177
$ws.Range("C19").Value = 0.021581
$ws.Range("D19").Value = 0.005915
$ws.Range("E19").Value = 0.006838
$ws.Range("F19").Value = 0.025339
$ws.Range("G19").Value = 0.023701
$ws.Range("H19").Value = 0.011036
$ws.Range("I19").Value = 0.014349

$ws.Range("B22").Value = 0.025178
$ws.Range("C22").Value = 0.01567
$ws.Range("D22").Value = 0.005819
$ws.Range("E22").Value = 0.003905
$ws.Range("F22").Value = 0.052021
$ws.Range("G22").Value = 0.429414
$ws.Range("H22").Value = 0.013547
$ws.Range("I22").Value = 0.011056

$ws.Range("B23").Value = 0.035356
$ws.Range("C23").Value = 0.023622
$ws.Range("D23").Value = 0.006104
$ws.Range("E23").Value = 0.007022
$ws.Range("F23").Value = 0.02766
$ws.Range("G23").Value = 0.026011
$ws.Range("H23").Value = 0.011394
$ws.Range("I23").Value = 0.014496

$ws.Range("B31").Value = 0.004695
$ws.Range("C31").Value = 0.004284
$ws.Range("D31").Value = 0.000518
$ws.Range("E31").Value = 0.000966
$ws.Range("F31").Value = 0.005543
$ws.Range("G31").Value = 0.00677
$ws.Range("H31").Value = 0.01182
$ws.Range("I31").Value = 0.014

$ws.Range("B32").Value = 0.034606
$ws.Range("C32").Value = 0.021523
$ws.Range("D32").Value = 0.007467
$ws.Range("E32").Value = 0.008872
$ws.Range("F32").Value = 0.028989
$ws.Range("G32").Value = 0.025404
$ws.Range("H32").Value = 0.039563
$ws.Range("I32").Value = 0.042996

$ws.Range("B46").Value = 0.075967
$ws.Range("C46").Value = 0.052598
$ws.Range("D46").Value = 0.027946
$ws.Range("E46").Value = 0.004708
$ws.Range("F46").Value = 0.037857
$ws.Range("G46").Value = 0.052748
$ws.Range("H46").Value = 0.021683

$ws.Range("B47").Value = 0.049123
$ws.Range("C47").Value = 0.031572
$ws.Range("D47").Value = 0.013261
$ws.Range("E47").Value = 0.010089
$ws.Range("F47").Value = 0.041035
$ws.Range("G47").Value = 0.04056
$ws.Range("H47").Value = 0.027946

# Narrow column I (previously widened for long N/A / skydome-related text)
$ws.Columns.Item(9).ColumnWidth = 11

# Update the active selection/view position
$ws.Range("G41").Select()
